$wb = $excel.ActiveWorkbook

# Sheets 1-3 and 5: row 1 headers B1:E1 are year labels that get an "Ano " prefix.
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet 4: row 1 headers B1:E1 are year/interval labels that get an "Intervalo " prefix.
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo 2015"
$ws4.Range("C1").Value = "Intervalo 2015-2030"
$ws4.Range("D1").Value = "Intervalo 2031-2040"
$ws4.Range("E1").Value = "Intervalo 2041-2050"

# Sheet 6: row 1 header B1 only gets the "Ano " prefix.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano 2015"
